$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4330.699162622797
$ws.Range("C3").Value = 4109.194807013856
$ws.Range("C4").Value = 4040.991385841469
$ws.Range("C5").Value = 4040.991385841469
$ws.Range("C6").Value = 4040.991385841469
$ws.Range("C7").Value = 4009.472353879031
$ws.Range("C8").Value = 3902.286189674337
$ws.Range("C9").Value = 3902.286189674337
$ws.Range("C10").Value = 3888.194353691815
$ws.Range("C11").Value = 3888.194353691815
$ws.Range("C12").Value = 3888.194353691815
